$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily records (date serial, nuovi pos., somma mobile 7gg., somma mobile per 100mila ab.)
$data = @(
    @(44441, 7, 35, 86.89607229753216),
    @(44442, 11, 37, 91.86156214310542),
    @(44443, 4, 36, 89.37881722031878),
    @(44444, 1, 36, 89.37881722031878),
    @(44445, 4, 32, 79.44783752917226),
    @(44446, 5, 35, 86.89607229753216),
    @(44447, 2, 34, 84.41332737474552),
    @(44448, 2, 29, 71.99960276081235)
)

$lastRow = 366
$startRow = $lastRow + 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i

    # Copy formatting (style, incl. the date number format) from the last existing row
    # onto column A of the new row, mirroring how the sheet was originally extended.
    $ws.Range("A$lastRow").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)

    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
